# Remove footnote markers like " [1]" and collapse embedded line breaks
# into single spaces across the vaccine price-list workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Pediatric VFC Vaccine
$ws2 = $wb.Worksheets.Item(2)   # Adult Vaccine
$ws3 = $wb.Worksheets.Item(3)   # Pediatric Influenza Vaccine
$ws4 = $wb.Worksheets.Item(4)   # Adult Influenza Vaccine

# --- Sheet 1: Pediatric VFC Vaccine ---
$ws1.Range("A2").Value = "DTaP "
$ws1.Range("A3").Value = "DTaP "
$ws1.Range("A4").Value = "DTaP "
$ws1.Range("A5").Value = "DTaP-IPV "
$ws1.Range("A6").Value = "DTaP-IPV "
$ws1.Range("A7").Value = "DTaP-IPV "
$ws1.Range("A8").Value = "DTaP-Hep B-IPV "
$ws1.Range("A9").Value = "DTaP-IP-HI "
$ws1.Range("A10").Value = "e-IPV "
$ws1.Range("A11").Value = "Hepatitis A Pediatric "
$ws1.Range("A12").Value = "Hepatitis A Pediatric "
$ws1.Range("A13").Value = "Hepatitis A Pediatric "
$ws1.Range("A14").Value = "Hepatitis A-Hepatitis B 18 only "
$ws1.Range("A15").Value = "Hepatitis B  Pediatric/Adolescent"
$ws1.Range("A16").Value = "Hepatitis B  Pediatric/Adolescent"
$ws1.Range("B16").Value = "Recombivax HB"
$ws1.Range("A17").Value = "Hib "
$ws1.Range("A18").Value = "Hib "
$ws1.Range("A19").Value = "Hib "
$ws1.Range("A20").Value = "HPV - Human Papillomavirus 9-valent "
$ws1.Range("A21").Value = "MENB - Meningococcal Group B "
$ws1.Range("A22").Value = "MENB - Meningococcal Group B "
$ws1.Range("A23").Value = "MENB - Meningococcal Group B "
$ws1.Range("A24").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws1.Range("A25").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws1.Range("A26").Value = "Measles, Mumps and Rubella (MMR) "
$ws1.Range("A27").Value = "MMR/Varicella "
$ws1.Range("A28").Value = "Pneumococcal 13-valent  (Pediatric)"
$ws1.Range("A30").Value = "Rotavirus, Live, Oral, Pentavalent "
$ws1.Range("A31").Value = "Rotavirus, Live, Oral, Pentavalent "
$ws1.Range("A32").Value = "Rotavirus, Live, Oral, Oral "
$ws1.Range("A33").Value = "Tetanus and Diphtheria Toxoids "
$ws1.Range("A34").Value = "Tetanus and Diphtheria Toxoids "
$ws1.Range("A35").Value = "Tetanus and Diphtheria Toxoids "
$ws1.Range("A36").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A37").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A38").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A39").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws1.Range("A40").Value = "Varicella "

# --- Sheet 2: Adult Vaccine ---
$ws2.Range("A2").Value = "Hepatitis A-Adult "
$ws2.Range("A3").Value = "Hepatitis A-Adult "
$ws2.Range("A4").Value = "Hepatitis A Adult "
$ws2.Range("A5").Value = "Hepatitis A Adult "
$ws2.Range("A6").Value = "Hepatitis A-Hepatitis B Adult "
$ws2.Range("A7").Value = "Hepatitis B-Adult "
$ws2.Range("A8").Value = "Hepatitis B-Adult "
$ws2.Range("A9").Value = "HPV-Human Papillomavirus 9 Valent "
$ws2.Range("A10").Value = "Measles, Mumps,  Rubella-Adult "
$ws2.Range("A11").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws2.Range("A12").Value = "Meningococcal Conjugate (Groups A, C, Y and W-135) "
$ws2.Range("A13").Value = "MENB - Meningococcal Group B "
$ws2.Range("A14").Value = "MENB - Meningococcal Group B "
$ws2.Range("A15").Value = "MENB - Meningococcal Group B "
$ws2.Range("A16").Value = "Pneumococcal 13-valent  (Adult)"
$ws2.Range("A19").Value = "Tetanus and Diphtheria Toxoids "
$ws2.Range("A20").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A21").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A22").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A23").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "
$ws2.Range("A24").Value = "Varicella-Adult "

# --- Sheet 3: Pediatric Influenza Vaccine ---
$ws3.Range("A2").Value = "Influenza  (Age 6 months and older)"
$ws3.Range("B2").Value = "Fluzone Quadrivalent"
$ws3.Range("A3").Value = "Influenza  (Age 6-35 months)"
$ws3.Range("B3").Value = "Fluzone Quadrivalent Pediatric dose"
$ws3.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B4").Value = "Fluzone Quadrivalent"
$ws3.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B5").Value = "Fluzone Quadrivalent"
$ws3.Range("A6").Value = "Influenza  (Age 36 months and older)"
$ws3.Range("B6").Value = "Fluarix Quadrivalent"
$ws3.Range("A7").Value = "Influenza  (Age 6 months and older)"
$ws3.Range("B7").Value = "FluLaval Quadrivalent"
$ws3.Range("A8").Value = "Influenza  (Age 6 months and older)"
$ws3.Range("B8").Value = "FluLaval Quadrivalent"
$ws3.Range("A9").Value = "Influenza  (Age 4 years and older)"
$ws3.Range("A10").Value = "Influenza  (Age 4 years and older)"

# --- Sheet 4: Adult Influenza Vaccine ---
$ws4.Range("A2").Value = "Influenza  (Age 6 months and older)"
$ws4.Range("B2").Value = "Fluzone Quadrivalent"
$ws4.Range("A3").Value = "Influenza  (Age 36 months and older)"
$ws4.Range("B3").Value = "Fluzone Quadrivalent"
$ws4.Range("A4").Value = "Influenza  (Age 36 months and older)"
$ws4.Range("B4").Value = "Fluzone Quadrivalent"
$ws4.Range("A5").Value = "Influenza  (Age 36 months and older)"
$ws4.Range("B5").Value = "Fluarix Quadrivalent"
$ws4.Range("A6").Value = "Influenza  (Age 36 months and older)"
$ws4.Range("B6").Value = "FluLaval Quadrivalent"
$ws4.Range("A7").Value = "Influenza  (Age 9 years and older)"
$ws4.Range("D7").Value = "10 pack-1 dose syringe"
$ws4.Range("A8").Value = "Influenza  (Age 9 years and older)"
$ws4.Range("A9").Value = "Influenza  (Age 4 years and older)"
$ws4.Range("A10").Value = "Influenza  (Age 4 years and older)"
$ws4.Range("A11").Value = "Influenza  (Age 4 years and older)"
$ws4.Range("A12").Value = "Influenza  (Age 4 years and older)"
$ws4.Range("A13").Value = "Influenza  (Age 18 years and older)"
$ws4.Range("B13").Value = "Afluria Quadrivalent"
$ws4.Range("A14").Value = "Influenza  (Age 18 years and older)"
$ws4.Range("B14").Value = "Afluria Quadrivalent"
